# Applies the OOXML diff:
#  - Slide 5 ("C'est quoi LINUX ?"): reposition/resize 5 shapes (text boxes
#    and the picture) - a layout tweak from the design-ideas redo.
#  - Slide 13 ("A votre tour !"): split the single run into three runs
#    ("A " / "votre" / " tour !") so "votre" carries its own run properties
#    (spell-check flag).
#
# NOTE on the literal point values below: PowerPoint's COM object model
# stores Shape.Left/Top/Width/Height as points in a 32-bit (Single) float,
# which the host then converts back to EMU (1 pt = 12700 EMU) when writing
# the OOXML. A plain "EMU / 12700" division rounds to the *nearest*
# representable Single, which can land fractionally below the target EMU
# value and truncate down once converted back (e.g. 417.7250393700787 ->
# 5305107 instead of 5305108). The literals used here are instead the
# smallest points value that still round-trips (via Single) to exactly the
# target EMU offset/extent from the diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5: move/resize shapes
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# id=28 "ZoneTexte 2" -> off (5305108, 1883392) EMU, ext unchanged
$shZoneTexte2 = $s5.Shapes.Item(2)
$shZoneTexte2.Left = 417.7250671386719
$shZoneTexte2.Top  = 148.298583984375

# id=74 "Picture 4" -> off (934447, 3168761) EMU, ext (2635345, 3350028) EMU
$shPicture4 = $s5.Shapes.Item(3)
$shPicture4.Left   = 73.57850646972656
$shPicture4.Top    = 249.5087432861328
$shPicture4.Width  = 207.5074920654297
$shPicture4.Height = 263.78173828125

# id=154 "Titre 1" -> off (5299590, 733304) EMU, ext unchanged
$shTitre1 = $s5.Shapes.Item(5)
$shTitre1.Left = 417.2905578613281
$shTitre1.Top  = 57.740474700927734

# id=165 "ZoneTexte 164" -> off (5299590, 2986547) EMU, ext unchanged
$shZoneTexte164 = $s5.Shapes.Item(6)
$shZoneTexte164.Left = 417.2905578613281
$shZoneTexte164.Top  = 235.16119384765625

# id=167 "ZoneTexte 166" -> off (5299590, 4233829) EMU, ext unchanged
$shZoneTexte166 = $s5.Shapes.Item(7)
$shZoneTexte166.Left = 417.2905578613281
$shZoneTexte166.Top  = 333.37237548828125

# ---------------------------------------------------------------------------
# Slide 13: "A votre tour !" -> split into three runs: "A ", "votre", " tour !"
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shTitre = $s13.Shapes.Item(8)
$tr = $shTitre.TextFrame.TextRange

# Re-assert the full text (keeps it a single run / paragraph), then restyle
# each sub-range in turn - this is what makes PowerPoint split a run into
# several runs sharing the same visible formatting but distinct rPr nodes.
$tr.Text = "A votre tour !"

$run1 = $tr.Characters(1, 2)
$run1.Font.Size = 60

$run2 = $tr.Characters(3, 5)
$run2.Font.Size = 60

$run3 = $tr.Characters(8, 7)
$run3.Font.Size = 60
